# Updated symbol list on Fri Dec 23 14:06:52 UTC 2022 with GitHub Actions
#
# Refresh of the crypto-ranking snapshot: Price (col D), Volume(1h) label
# (col E) and Hora/hour (col G) are re-pulled for every ranked row, and a
# handful of rows were re-ordered by the provider, so their Coin (col B)
# and Link (col C) cells move too.
#
# NumberFormat is forced to text ("@") before each write so numeric-looking
# values (e.g. "245.99") are stored as text, matching the source file, and
# not re-typed into floating point numbers. Style is reset back to "Normal"
# afterwards so no stray text-format style lingers on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$cell = $ws.Cells.Item(2, 4)
$cell.NumberFormat = "@"
$cell.Value = '245.99'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(2, 7)
$cell.NumberFormat = "@"
$cell.Value = '14'
$cell.Style = "Normal"
# Row 3
$cell = $ws.Cells.Item(3, 4)
$cell.NumberFormat = "@"
$cell.Value = '22.11'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(3, 7)
$cell.NumberFormat = "@"
$cell.Value = '14'
$cell.Style = "Normal"
# Row 4
$cell = $ws.Cells.Item(4, 4)
$cell.NumberFormat = "@"
$cell.Value = '5.423'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(4, 7)
$cell.NumberFormat = "@"
$cell.Value = '14'
$cell.Style = "Normal"
# Row 5
$cell = $ws.Cells.Item(5, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.05856'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(5, 7)
$cell.NumberFormat = "@"
$cell.Value = '14'
$cell.Style = "Normal"
# Row 6
$cell = $ws.Cells.Item(6, 4)
$cell.NumberFormat = "@"
$cell.Value = '3.394'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(6, 7)
$cell.NumberFormat = "@"
$cell.Value = '14'
$cell.Style = "Normal"
# Row 7
$cell = $ws.Cells.Item(7, 4)
$cell.NumberFormat = "@"
$cell.Value = '6.361'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(7, 7)
$cell.NumberFormat = "@"
$cell.Value = '14'
$cell.Style = "Normal"
# Row 8
$cell = $ws.Cells.Item(8, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.8150'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(8, 7)
$cell.NumberFormat = "@"
$cell.Value = '14'
$cell.Style = "Normal"
# Row 9
$cell = $ws.Cells.Item(9, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.002'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(9, 5)
$cell.NumberFormat = "@"
$cell.Value = '8FTXTokenFTTBestin24h'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(9, 7)
$cell.NumberFormat = "@"
$cell.Value = '14'
$cell.Style = "Normal"
# Row 10
$cell = $ws.Cells.Item(10, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.1430'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(10, 7)
$cell.NumberFormat = "@"
$cell.Value = '14'
$cell.Style = "Normal"
# Row 11
$cell = $ws.Cells.Item(11, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.07439'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(11, 7)
$cell.NumberFormat = "@"
$cell.Value = '14'
$cell.Style = "Normal"
# Row 12
$cell = $ws.Cells.Item(12, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.03433'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(12, 7)
$cell.NumberFormat = "@"
$cell.Value = '14'
$cell.Style = "Normal"
# Row 13
$cell = $ws.Cells.Item(13, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.03023'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(13, 7)
$cell.NumberFormat = "@"
$cell.Value = '14'
$cell.Style = "Normal"
# Row 14
$cell = $ws.Cells.Item(14, 4)
$cell.NumberFormat = "@"
$cell.Value = '4.186'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(14, 7)
$cell.NumberFormat = "@"
$cell.Value = '14'
$cell.Style = "Normal"
# Row 15
$cell = $ws.Cells.Item(15, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.09403'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(15, 7)
$cell.NumberFormat = "@"
$cell.Value = '14'
$cell.Style = "Normal"
# Row 16
$cell = $ws.Cells.Item(16, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.001600'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(16, 7)
$cell.NumberFormat = "@"
$cell.Value = '14'
$cell.Style = "Normal"
# Row 17
$cell = $ws.Cells.Item(17, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.04824'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(17, 7)
$cell.NumberFormat = "@"
$cell.Value = '14'
$cell.Style = "Normal"
# Row 18
$cell = $ws.Cells.Item(18, 2)
$cell.NumberFormat = "@"
$cell.Value = 'One'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(18, 3)
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(18, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.0005894'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(18, 5)
$cell.NumberFormat = "@"
$cell.Value = '17OneONEWorstin24h'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(18, 7)
$cell.NumberFormat = "@"
$cell.Value = '14'
$cell.Style = "Normal"
# Row 19
$cell = $ws.Cells.Item(19, 2)
$cell.NumberFormat = "@"
$cell.Value = 'TigerCash'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(19, 3)
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(19, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.006018'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(19, 5)
$cell.NumberFormat = "@"
$cell.Value = '18TigerCashTCH'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(19, 7)
$cell.NumberFormat = "@"
$cell.Value = '14'
$cell.Style = "Normal"
# Row 20
$cell = $ws.Cells.Item(20, 2)
$cell.NumberFormat = "@"
$cell.Value = 'HotbitToken'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(20, 3)
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(20, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.004098'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(20, 5)
$cell.NumberFormat = "@"
$cell.Value = '19HotbitTokenHTB'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(20, 7)
$cell.NumberFormat = "@"
$cell.Value = '14'
$cell.Style = "Normal"
# Row 21
$cell = $ws.Cells.Item(21, 2)
$cell.NumberFormat = "@"
$cell.Value = 'BitKan'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(21, 3)
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(21, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.0009980'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(21, 5)
$cell.NumberFormat = "@"
$cell.Value = '20BitKanKAN'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(21, 7)
$cell.NumberFormat = "@"
$cell.Value = '14'
$cell.Style = "Normal"
# Row 22
$cell = $ws.Cells.Item(22, 2)
$cell.NumberFormat = "@"
$cell.Value = 'NitroEx'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(22, 3)
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(22, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.0001502'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(22, 5)
$cell.NumberFormat = "@"
$cell.Value = '21NitroExNTX'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(22, 7)
$cell.NumberFormat = "@"
$cell.Value = '14'
$cell.Style = "Normal"
# Row 23
$cell = $ws.Cells.Item(23, 2)
$cell.NumberFormat = "@"
$cell.Value = 'LEO'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(23, 3)
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(23, 4)
$cell.NumberFormat = "@"
$cell.Value = '3.698'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(23, 5)
$cell.NumberFormat = "@"
$cell.Value = '22LEOLEO'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(23, 7)
$cell.NumberFormat = "@"
$cell.Value = '14'
$cell.Style = "Normal"
# Row 24
$cell = $ws.Cells.Item(24, 2)
$cell.NumberFormat = "@"
$cell.Value = 'BTSEToken'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(24, 3)
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(24, 4)
$cell.NumberFormat = "@"
$cell.Value = '2.217'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(24, 5)
$cell.NumberFormat = "@"
$cell.Value = '23BTSETokenBTSE'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(24, 7)
$cell.NumberFormat = "@"
$cell.Value = '14'
$cell.Style = "Normal"
# Row 25
$cell = $ws.Cells.Item(25, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.3239'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(25, 7)
$cell.NumberFormat = "@"
$cell.Value = '14'
$cell.Style = "Normal"
# Row 26
$cell = $ws.Cells.Item(26, 7)
$cell.NumberFormat = "@"
$cell.Value = '14'
$cell.Style = "Normal"
# Row 27
$cell = $ws.Cells.Item(27, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.0001291'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(27, 5)
$cell.NumberFormat = "@"
$cell.Value = '26UpBotsUBXT'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(27, 7)
$cell.NumberFormat = "@"
$cell.Value = '14'
$cell.Style = "Normal"
# Row 28
$cell = $ws.Cells.Item(28, 7)
$cell.NumberFormat = "@"
$cell.Value = '14'
$cell.Style = "Normal"
# Row 29
$cell = $ws.Cells.Item(29, 7)
$cell.NumberFormat = "@"
$cell.Value = '14'
$cell.Style = "Normal"
# Row 30
$cell = $ws.Cells.Item(30, 7)
$cell.NumberFormat = "@"
$cell.Value = '14'
$cell.Style = "Normal"
# Row 31
$cell = $ws.Cells.Item(31, 7)
$cell.NumberFormat = "@"
$cell.Value = '14'
$cell.Style = "Normal"
# Row 32
$cell = $ws.Cells.Item(32, 7)
$cell.NumberFormat = "@"
$cell.Value = '14'
$cell.Style = "Normal"
# Row 33
$cell = $ws.Cells.Item(33, 7)
$cell.NumberFormat = "@"
$cell.Value = '14'
$cell.Style = "Normal"
# Row 34
$cell = $ws.Cells.Item(34, 7)
$cell.NumberFormat = "@"
$cell.Value = '14'
$cell.Style = "Normal"
# Row 35
$cell = $ws.Cells.Item(35, 7)
$cell.NumberFormat = "@"
$cell.Value = '14'
$cell.Style = "Normal"
# Row 36
$cell = $ws.Cells.Item(36, 7)
$cell.NumberFormat = "@"
$cell.Value = '14'
$cell.Style = "Normal"
# Row 37
$cell = $ws.Cells.Item(37, 7)
$cell.NumberFormat = "@"
$cell.Value = '14'
$cell.Style = "Normal"
# Row 38
$cell = $ws.Cells.Item(38, 7)
$cell.NumberFormat = "@"
$cell.Value = '14'
$cell.Style = "Normal"
# Row 39
$cell = $ws.Cells.Item(39, 7)
$cell.NumberFormat = "@"
$cell.Value = '14'
$cell.Style = "Normal"
# Row 40
$cell = $ws.Cells.Item(40, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.03861'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(40, 7)
$cell.NumberFormat = "@"
$cell.Value = '14'
$cell.Style = "Normal"
# Row 41
$cell = $ws.Cells.Item(41, 2)
$cell.NumberFormat = "@"
$cell.Value = 'BKEXToken'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(41, 3)
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(41, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.1076'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(41, 5)
$cell.NumberFormat = "@"
$cell.Value = '40BKEXTokenBKK'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(41, 7)
$cell.NumberFormat = "@"
$cell.Value = '14'
$cell.Style = "Normal"
# Row 42
$cell = $ws.Cells.Item(42, 2)
$cell.NumberFormat = "@"
$cell.Value = 'CEJI'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(42, 3)
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(42, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.002413'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(42, 5)
$cell.NumberFormat = "@"
$cell.Value = '41CEJICEJI'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(42, 7)
$cell.NumberFormat = "@"
$cell.Value = '14'
$cell.Style = "Normal"
# Row 43
$cell = $ws.Cells.Item(43, 2)
$cell.NumberFormat = "@"
$cell.Value = 'KickToken'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(43, 3)
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(43, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.003023'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(43, 5)
$cell.NumberFormat = "@"
$cell.Value = '42KickTokenKICK'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(43, 7)
$cell.NumberFormat = "@"
$cell.Value = '14'
$cell.Style = "Normal"
# Row 44
$cell = $ws.Cells.Item(44, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.006250'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(44, 7)
$cell.NumberFormat = "@"
$cell.Value = '14'
$cell.Style = "Normal"
# Row 45
$cell = $ws.Cells.Item(45, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.00005626'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(45, 7)
$cell.NumberFormat = "@"
$cell.Value = '14'
$cell.Style = "Normal"
# Row 46
$cell = $ws.Cells.Item(46, 7)
$cell.NumberFormat = "@"
$cell.Value = '14'
$cell.Style = "Normal"
# Row 47
$cell = $ws.Cells.Item(47, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.4002'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(47, 7)
$cell.NumberFormat = "@"
$cell.Value = '14'
$cell.Style = "Normal"
# Row 48
$cell = $ws.Cells.Item(48, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.1420'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(48, 7)
$cell.NumberFormat = "@"
$cell.Value = '14'
$cell.Style = "Normal"
# Row 49
$cell = $ws.Cells.Item(49, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.00002101'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(49, 7)
$cell.NumberFormat = "@"
$cell.Value = '14'
$cell.Style = "Normal"
# Row 50
$cell = $ws.Cells.Item(50, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.01010'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(50, 7)
$cell.NumberFormat = "@"
$cell.Value = '14'
$cell.Style = "Normal"
# Row 51
$cell = $ws.Cells.Item(51, 7)
$cell.NumberFormat = "@"
$cell.Value = '14'
$cell.Style = "Normal"
